$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 9 / Row 10 swap: Dogecoin and LidoStakedEther exchange rank positions ---
$ws.Range("B9").Value = "Dogecoin"
$ws.Range("C9").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D9").Value = "'0.172"
$ws.Range("E9").Value = "  -1.34%  "

$ws.Range("B10").Value = "LidoStakedEther"
$ws.Range("C10").Value = "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth"
$ws.Range("D10").Value = "'2.637.09"
$ws.Range("E10").Value = "  -0.89%  "

# --- Remaining price / volume updates ---
$ws.Range("D2").Value = "'72.271.49"
$ws.Range("E2").Value = "  -0.16%  "
$ws.Range("D3").Value = "'2.637.78"
$ws.Range("E3").Value = "  -0.90%  "
$ws.Range("D5").Value = "'589.07"
$ws.Range("E5").Value = "  -2.20%  "
$ws.Range("D6").Value = "'173.72"
$ws.Range("E6").Value = "  -2.83%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").Value = "'0.518"
$ws.Range("E8").Value = "  -1.25%  "
$ws.Range("E11").Value = "  +1.10%  "
$ws.Range("D12").Value = "'0.354"
$ws.Range("E12").Value = "  -1.03%  "
$ws.Range("D13").Value = "'4.93"
$ws.Range("E13").Value = "  -1.78%  "
$ws.Range("D14").Value = "'3.131.98"
$ws.Range("E14").Value = "  -0.47%  "
$ws.Range("D15").Value = "'0.0000185"
$ws.Range("E15").Value = "  -2.26%  "
$ws.Range("D16").Value = "'72.192.37"
$ws.Range("E16").Value = "  -0.11%  "
$ws.Range("D17").Value = "'25.83"
$ws.Range("E17").Value = "  -2.74%  "
$ws.Range("D18").Value = "'2.640.11"
$ws.Range("E18").Value = "  -1.03%  "
$ws.Range("D19").Value = "'12.12"
$ws.Range("E19").Value = "  +1.05%  "
$ws.Range("D20").Value = "'7.95"
$ws.Range("E20").Value = "  -0.68%  "
$ws.Range("D21").Value = "'372.13"
$ws.Range("E21").Value = "  -1.77%  "
$ws.Range("D22").Value = "'4.14"
$ws.Range("E22").Value = "  -1.59%  "
$ws.Range("D23").Value = "'2.04"
$ws.Range("E23").Value = "  -2.14%  "
$ws.Range("E24").Value = "  -0.02%  "
$ws.Range("D25").Value = "'70.91"
$ws.Range("E25").Value = "  -2.47%  "
$ws.Range("D26").Value = "'4.24"
$ws.Range("E26").Value = "  -3.91%  "
$ws.Range("D27").Value = "'9.62"
$ws.Range("E27").Value = "  -3.78%  "
$ws.Range("D28").Value = "'2.772.82"
$ws.Range("E28").Value = "  -1.22%  "
$ws.Range("D29").Value = "'0.999"
$ws.Range("E29").Value = "  +0.17%  "
$ws.Range("D30").Value = "'0.0₃0952"
$ws.Range("E30").Value = "  -0.20%  "
$ws.Range("D31").Value = "'7.97"
$ws.Range("E31").Value = "  -3.35%  "
$ws.Range("D32").Value = "'495.04"
$ws.Range("E32").Value = "  -5.13%  "
$ws.Range("E33").Value = "  -2.80%  "
$ws.Range("E34").Value = "  -1.23%  "
$ws.Range("D36").Value = "'161.46"
$ws.Range("E36").Value = "  -1.12%  "
$ws.Range("D37").Value = "'19.25"
$ws.Range("E37").Value = "  -1.31%  "
$ws.Range("D38").Value = "'0.112"
$ws.Range("E38").Value = "  +0.40%  "
$ws.Range("D39").Value = "'18.87"
$ws.Range("E39").Value = "  -1.31%  "
$ws.Range("E40").Value = "  -3.03%  "
$ws.Range("E41").Value = "  -0.09%  "
$ws.Range("E42").Value = "  -6.00%  "
$ws.Range("D43").Value = "'2.55"
$ws.Range("E43").Value = "  -2.31%  "
$ws.Range("D44").Value = "'4.89"
$ws.Range("E44").Value = "  -3.68%  "
$ws.Range("D45").Value = "'0.326"
$ws.Range("E45").Value = "  -2.82%  "
$ws.Range("D46").Value = "'39.06"
$ws.Range("E46").Value = "  -0.79%  "
$ws.Range("D47").Value = "'152.53"
$ws.Range("E47").Value = "  -0.40%  "
$ws.Range("D48").Value = "'3.64"
$ws.Range("E48").Value = "  -2.72%  "
$ws.Range("D49").Value = "'0.545"
$ws.Range("E49").Value = "  -0.54%  "
$ws.Range("D50").Value = "'1.67"
$ws.Range("E50").Value = "  -2.60%  "
$ws.Range("E51").Value = "  -1.89%  "
